# Auto-generated edit script: update scraped horarios data for Linea 141
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2,1).Value = 'Última actualización: 08:39:44'
$ws.Cells.Item(3,1).Value = 'Total filas: 115'
$ws.Cells.Item(47,1).Value = '05:18:56'
$ws.Cells.Item(47,3).Value = '15_ABASTO'
$ws.Cells.Item(47,4).Value = 106
$ws.Cells.Item(48,1).Value = '05:49:40'
$ws.Cells.Item(48,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(48,4).Value = 75
$ws.Cells.Item(62,1).Value = '06:15:04'
$ws.Cells.Item(62,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(62,4).Value = 77
$ws.Cells.Item(63,1).Value = '05:49:40'
$ws.Cells.Item(63,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(63,4).Value = 103
$ws.Cells.Item(98,1).Value = '08:39:44'
$ws.Cells.Item(98,2).Value = '09:04'
$ws.Cells.Item(98,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(98,4).Value = 25
$ws.Cells.Item(99,1).Value = '08:21:50'
$ws.Cells.Item(99,2).Value = '09:07'
$ws.Cells.Item(99,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(99,4).Value = 46
$ws.Cells.Item(100,1).Value = '07:20:40'
$ws.Cells.Item(100,2).Value = '09:10'
$ws.Cells.Item(100,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(100,4).Value = 110
$ws.Cells.Item(101,1).Value = '08:21:50'
$ws.Cells.Item(101,2).Value = '09:13'
$ws.Cells.Item(101,3).Value = '10_OLMOS'
$ws.Cells.Item(101,4).Value = 52
$ws.Cells.Item(102,1).Value = '07:20:40'
$ws.Cells.Item(102,2).Value = '09:16'
$ws.Cells.Item(102,3).Value = '27_EL RETIRO'
$ws.Cells.Item(102,4).Value = 116
$ws.Cells.Item(103,1).Value = '08:21:50'
$ws.Cells.Item(103,2).Value = '09:21'
$ws.Cells.Item(103,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(103,4).Value = 60
$ws.Cells.Item(104,1).Value = '07:59:28'
$ws.Cells.Item(104,3).Value = '16_SANTA ANA'
$ws.Cells.Item(104,4).Value = 83
$ws.Cells.Item(105,2).Value = '09:22'
$ws.Cells.Item(105,3).Value = '17_ROMERO'
$ws.Cells.Item(105,4).Value = 95
$ws.Cells.Item(106,1).Value = '07:47:32'
$ws.Cells.Item(106,2).Value = '09:23'
$ws.Cells.Item(106,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(106,4).Value = 96
$ws.Cells.Item(107,1).Value = '08:21:50'
$ws.Cells.Item(107,2).Value = '09:29'
$ws.Cells.Item(107,3).Value = '16_SANTA ANA'
$ws.Cells.Item(107,4).Value = 68
$ws.Cells.Item(108,2).Value = '09:32'
$ws.Cells.Item(108,3).Value = '15_ABASTO'
$ws.Cells.Item(108,4).Value = 105
$ws.Cells.Item(109,1).Value = '07:47:32'
$ws.Cells.Item(109,2).Value = '09:33'
$ws.Cells.Item(109,3).Value = '10_OLMOS'
$ws.Cells.Item(109,4).Value = 106
$ws.Cells.Item(110,1).Value = '08:39:44'
$ws.Cells.Item(110,2).Value = '09:34'
$ws.Cells.Item(110,3).Value = '16_SANTA ANA'
$ws.Cells.Item(110,4).Value = 55
$ws.Cells.Item(111,1).Value = '08:39:44'
$ws.Cells.Item(111,2).Value = '09:34'
$ws.Cells.Item(111,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(111,4).Value = 55
$ws.Cells.Item(112,2).Value = '09:41'
$ws.Cells.Item(112,3).Value = '215C_EL PATO'
$ws.Cells.Item(112,4).Value = 80
$ws.Cells.Item(113,1).Value = '07:47:32'
$ws.Cells.Item(113,2).Value = '09:42'
$ws.Cells.Item(113,3).Value = '215C_EL PATO'
$ws.Cells.Item(113,4).Value = 115
$ws.Cells.Item(114,1).Value = '07:47:32'
$ws.Cells.Item(114,2).Value = '09:43'
$ws.Cells.Item(114,3).Value = '14_ABASTO'
$ws.Cells.Item(114,4).Value = 116
$ws.Cells.Item(114,5).Value = 'LP1912'
$ws.Cells.Item(115,1).Value = '08:39:44'
$ws.Cells.Item(115,2).Value = '10:06'
$ws.Cells.Item(115,3).Value = '10_OLMOS'
$ws.Cells.Item(115,4).Value = 87
$ws.Cells.Item(115,5).Value = 'LP1912'
$ws.Cells.Item(116,1).Value = '08:21:50'
$ws.Cells.Item(116,2).Value = '10:10'
$ws.Cells.Item(116,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(116,4).Value = 109
$ws.Cells.Item(116,5).Value = 'LP1912'
$ws.Cells.Item(117,1).Value = '08:21:50'
$ws.Cells.Item(117,2).Value = '10:12'
$ws.Cells.Item(117,3).Value = '15_ABASTO'
$ws.Cells.Item(117,4).Value = 111
$ws.Cells.Item(117,5).Value = 'LP1912'
$ws.Cells.Item(118,1).Value = '08:39:44'
$ws.Cells.Item(118,2).Value = '10:21'
$ws.Cells.Item(118,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(118,4).Value = 102
$ws.Cells.Item(118,5).Value = 'LP1912'
$ws.Cells.Item(119,1).Value = '08:39:44'
$ws.Cells.Item(119,2).Value = '10:22'
$ws.Cells.Item(119,3).Value = '17_ROMERO'
$ws.Cells.Item(119,4).Value = 103
$ws.Cells.Item(119,5).Value = 'LP1912'
$ws.Cells.Item(120,1).Value = '08:39:44'
$ws.Cells.Item(120,2).Value = '10:26'
$ws.Cells.Item(120,3).Value = '215A_EL PATO'
$ws.Cells.Item(120,4).Value = 107
$ws.Cells.Item(120,5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2,1).Value = 'Última actualización: 08:39:44'
$ws.Cells.Item(3,1).Value = 'Total filas: 17'
$ws.Cells.Item(22,1).Value = '08:39:44'
$ws.Cells.Item(22,2).Value = '10:26'
$ws.Cells.Item(22,3).Value = '215A_EL PATO'
$ws.Cells.Item(22,4).Value = 107
$ws.Cells.Item(22,5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2,1).Value = 'Última actualización: 08:39:44'
$ws.Cells.Item(3,1).Value = 'Total filas: 26'
$ws.Cells.Item(28,1).Value = '08:39:44'
$ws.Cells.Item(28,2).Value = '08:44'
$ws.Cells.Item(28,3).Value = '215A_LA PLATA'
$ws.Cells.Item(28,4).Value = 5
$ws.Cells.Item(28,5).Value = 'L6173'
$ws.Cells.Item(29,1).Value = '07:20:40'
$ws.Cells.Item(29,2).Value = '09:08'
$ws.Cells.Item(29,4).Value = 108
$ws.Cells.Item(30,1).Value = '07:47:32'
$ws.Cells.Item(30,2).Value = '09:09'
$ws.Cells.Item(30,3).Value = '215D_LA PLATA'
$ws.Cells.Item(30,4).Value = 82
$ws.Cells.Item(30,5).Value = 'L6203'
$ws.Cells.Item(31,1).Value = '08:21:50'
$ws.Cells.Item(31,2).Value = '10:02'
$ws.Cells.Item(31,3).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(31,4).Value = 101
$ws.Cells.Item(31,5).Value = 'L6173'

